# Auto-generated edit script: update crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.849.99'
$ws.Cells.Item(2, 5).Value = '  +0.14%  '

$ws.Cells.Item(3, 4).Value = '1.635.51'
$ws.Cells.Item(3, 5).Value = '  +0.46%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '215.25'
$ws.Cells.Item(5, 5).Value = '  -0.02%  '

$ws.Cells.Item(7, 5).Value = '  -0.05%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.258'
$ws.Cells.Item(8, 5).Value = '  -0.19%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.0642'
$ws.Cells.Item(9, 5).Value = '  +0.31%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '20.10'
$ws.Cells.Item(10, 5).Value = '  +3.81%  '

$ws.Cells.Item(11, 5).Value = '  +0.28%  '

$ws.Cells.Item(12, 4).Value = '1.666.23'
$ws.Cells.Item(12, 5).Value = '  +2.27%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '4.25'
$ws.Cells.Item(13, 5).Value = '  -0.03%  '

$ws.Cells.Item(14, 4).Value = '1.862.29'
$ws.Cells.Item(14, 5).Value = '  +0.55%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.562'
$ws.Cells.Item(15, 5).Value = '  +0.72%  '

$ws.Cells.Item(16, 4).Value = '0.0₃0764'
$ws.Cells.Item(16, 5).Value = '  +1.39%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '63.20'
$ws.Cells.Item(17, 5).Value = '  -0.49%  '

$ws.Cells.Item(18, 4).Value = '25.857.54'
$ws.Cells.Item(18, 5).Value = '  +0.18%  '

$ws.Cells.Item(19, 5).Value = '  -0.09%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '194.37'
$ws.Cells.Item(20, 5).Value = '  +0.16%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.37'

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '9.91'
$ws.Cells.Item(22, 5).Value = '  +1.08%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '6.18'
$ws.Cells.Item(23, 5).Value = '  +3.03%  '

$ws.Cells.Item(24, 5).Value = '  -0.04%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.74'
$ws.Cells.Item(25, 5).Value = '  -3.73%  '

$ws.Cells.Item(26, 5).Value = '  -2.05%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.123'
$ws.Cells.Item(27, 5).Value = '  -5.49%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '6.83'
$ws.Cells.Item(28, 5).Value = '  +1.21%  '

$ws.Cells.Item(29, 5).Value = '  +0.78%  '

$ws.Cells.Item(30, 5).Value = '  +0.19%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.0493'
$ws.Cells.Item(31, 5).Value = '  +1.12%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.30'
$ws.Cells.Item(32, 5).Value = '  +0.14%  '

$ws.Cells.Item(33, 5).Value = '  +1.59%  '

$ws.Cells.Item(34, 5).Value = '  +0.73%  '

$ws.Cells.Item(35, 5).Value = '  +0.62%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.903'
$ws.Cells.Item(36, 5).Value = '  +0.85%  '

$ws.Cells.Item(37, 5).Value = '  +1.38%  '

$ws.Cells.Item(38, 2).Value = 'Maker'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(38, 4).Value = '1.119.85'
$ws.Cells.Item(38, 5).Value = '  -1.52%  '

$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.547'
$ws.Cells.Item(39, 5).Value = '  -0.28%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.0157'
$ws.Cells.Item(40, 5).Value = '  +0.76%  '

$ws.Cells.Item(41, 5).Value = '  -0.30%  '

$ws.Cells.Item(42, 5).Value = '  -1.82%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '99.36'
$ws.Cells.Item(43, 5).Value = '  +2.05%  '

$ws.Cells.Item(44, 5).Value = '  +0.70%  '

$ws.Cells.Item(45, 5).Value = '  +1.14%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '55.37'
$ws.Cells.Item(46, 5).Value = '  +0.87%  '

$ws.Cells.Item(47, 5).Value = '  -4.37%  '

$ws.Cells.Item(48, 5).Value = '  -0.45%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '7.63'
$ws.Cells.Item(49, 5).Value = '  +0.98%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.999'
$ws.Cells.Item(50, 5).Value = '  -0.24%  '

$ws.Cells.Item(51, 5).Value = '  -0.06%  '
